$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2 = @("30.055.24", "  +4.06%  ")
    3 = @("1.896.05", "  +4.35%  ")
    4 = @("1.001", "  +0.18%  ")
    5 = @("249.02", "  +1.03%  ")
    6 = @("1.001", "  +0.10%  ")
    7 = @("0.4977", "  +0.84%  ")
    8 = @("45.03", "  +3.16%  ")
    9 = @("0.2964", "  +6.80%  ")
    10 = @("0.06651", "  +4.04%  ")
    11 = @("1.898.92", "  +4.54%  ")
    12 = @("17.03", "  +1.79%  ")
    13 = @("0.07249", "  +2.55%  ")
    14 = @("0.6792", "  +5.64%  ")
    15 = @("86.33", "  +2.76%  ")
    16 = @("4.874", "  +4.38%  ")
    17 = @("30.059.05", "  +3.98%  ")
    18 = @("0.000007987", "  +9.39%  ")
    19 = @("1.001", "  +0.12%  ")
    20 = @("12.92", "  +5.86%  ")
    21 = @("2.147.46", "  +5.21%  ")
    22 = @("1.002", "  +0.27%  ")
    23 = @("4.775", "  +4.87%  ")
    24 = @("5.682", "  +6.16%  ")
    25 = @("9.264", "  +5.26%  ")
    26 = @("148.20", "  +2.31%  ")
    27 = @("132.08", "  +2.31%  ")
    28 = @("16.81", "  +2.61%  ")
    29 = @("1.969", "  +4.76%  ")
    30 = @("1.384", "  -1.41%  ")
    31 = @("4.232", "  +2.58%  ")
    32 = @("0.08757", "  +4.93%  ")
    33 = @("3.944", "  +4.36%  ")
    34 = @("0.05101", "  +3.88%  ")
    35 = @("1.129", "  +3.00%  ")
    36 = @("0.7034", "  +4.50%  ")
    37 = @("2.687", "  -0.27%  ")
    38 = @("2.788", "  +3.02%  ")
    39 = @("2.232", "  -2.55%  ")
    40 = @("0.9535", "  +0.84%  ")
    41 = @("0.01667", "  +5.26%  ")
    42 = @("5.998", "  -2.57%  ")
    43 = @("0.4244", "  +3.96%  ")
    44 = @("0.9997", "  +0.01%  ")
    45 = @("103.37", "  +2.84%  ")
    46 = @("7.490", "  +4.54%  ")
    47 = @("0.1265", "  +3.44%  ")
    48 = @("0.05750", "  +4.24%  ")
    49 = @("32.93", "  +4.15%  ")
    50 = @("8.304", "  +2.70%  ")
    51 = @("0.3749", "  +3.78%  ")
}

foreach ($row in ($updates.Keys | Sort-Object { [int]$_ })) {
    $vals = $updates[$row]
    $dVal = $vals[0]
    $eVal = $vals[1]

    # Columns D hold text that sometimes looks like a plain number (e.g. "1.001").
    # Force those to stay text (matching the source inlineStr cells) instead of
    # letting Excel auto-convert them to a numeric value.
    $dCell = $ws.Cells.Item($row, 4)
    $looksNumeric = $dVal -match "^[+-]?[0-9]*\.?[0-9]+([eE][+-]?[0-9]+)?$"
    if ($looksNumeric) {
        $dCell.NumberFormat = "@"
        $dCell.Value = $dVal
        $dCell.Style = "Normal"
    } else {
        $dCell.Value = $dVal
    }

    $eCell = $ws.Cells.Item($row, 5)
    $eCell.Value = $eVal
}
